$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 19:18"

# Update country rows: names (re-sorted by Casos totales) and refreshed case counts
$ws.Cells.Item(4,1).Value = "Estados Unidos"
$ws.Cells.Item(4,2).Value = 2959188
$ws.Cells.Item(4,3).Value = 23418
$ws.Cells.Item(4,4).Value = 1261420
$ws.Cells.Item(4,5).Value = 1565350
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 100
$ws.Cells.Item(4,8).Value = 132418

$ws.Cells.Item(5,1).Value = "Brasil"
$ws.Cells.Item(5,2).Value = 1579837
$ws.Cells.Item(5,3).Value = 1461
$ws.Cells.Item(5,4).Value = 978615
$ws.Cells.Item(5,5).Value = 536839
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 18
$ws.Cells.Item(5,8).Value = 64383

$ws.Cells.Item(6,1).Value = "India"
$ws.Cells.Item(6,2).Value = 697069
$ws.Cells.Item(6,3).Value = 23165
$ws.Cells.Item(6,4).Value = 424885
$ws.Cells.Item(6,5).Value = 252485
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(6,7).Value = 420
$ws.Cells.Item(6,8).Value = 19699

$ws.Cells.Item(17,1).Value = "Turquia"
$ws.Cells.Item(17,2).Value = 205758
$ws.Cells.Item(17,3).Value = 1148
$ws.Cells.Item(17,4).Value = 180680
$ws.Cells.Item(17,5).Value = 19853
$ws.Cells.Item(17,6).Value = 0
$ws.Cells.Item(17,7).Value = 19
$ws.Cells.Item(17,8).Value = 5225

$ws.Cells.Item(18,1).Value = "Alemania"
$ws.Cells.Item(18,2).Value = 197460
$ws.Cells.Item(18,3).Value = 42
$ws.Cells.Item(18,4).Value = 181700
$ws.Cells.Item(18,5).Value = 6675
$ws.Cells.Item(18,6).Value = 0
$ws.Cells.Item(18,7).Value = 4
$ws.Cells.Item(18,8).Value = 9085

$ws.Cells.Item(40,1).Value = "Singapur"
$ws.Cells.Item(40,2).Value = 44800
$ws.Cells.Item(40,3).Value = 136
$ws.Cells.Item(40,4).Value = 40441
$ws.Cells.Item(40,5).Value = 4333
$ws.Cells.Item(40,6).Value = 0
$ws.Cells.Item(40,7).Value = 0
$ws.Cells.Item(40,8).Value = 26

$ws.Cells.Item(44,1).Value = "Republica Dominicana"
$ws.Cells.Item(44,2).Value = 37425
$ws.Cells.Item(44,3).Value = 1241
$ws.Cells.Item(44,4).Value = 18943
$ws.Cells.Item(44,5).Value = 17688
$ws.Cells.Item(44,6).Value = 0
$ws.Cells.Item(44,7).Value = 8
$ws.Cells.Item(44,8).Value = 794

$ws.Cells.Item(45,1).Value = "Panama"
$ws.Cells.Item(45,2).Value = 36983
$ws.Cells.Item(45,3).Value = 0
$ws.Cells.Item(45,4).Value = 17761
$ws.Cells.Item(45,5).Value = 18502
$ws.Cells.Item(45,6).Value = 0
$ws.Cells.Item(45,7).Value = 0
$ws.Cells.Item(45,8).Value = 720

$ws.Cells.Item(49,1).Value = "Israel"
$ws.Cells.Item(49,2).Value = 29787
$ws.Cells.Item(49,3).Value = 617
$ws.Cells.Item(49,4).Value = 17916
$ws.Cells.Item(49,5).Value = 11540
$ws.Cells.Item(49,6).Value = 0
$ws.Cells.Item(49,7).Value = 1
$ws.Cells.Item(49,8).Value = 331

$ws.Cells.Item(51,1).Value = "Barein"
$ws.Cells.Item(51,2).Value = 28857
$ws.Cells.Item(51,3).Value = 0
$ws.Cells.Item(51,4).Value = 23959
$ws.Cells.Item(51,5).Value = 4801
$ws.Cells.Item(51,6).Value = 0
$ws.Cells.Item(51,7).Value = 1
$ws.Cells.Item(51,8).Value = 97

$ws.Cells.Item(54,1).Value = "Irlanda"
$ws.Cells.Item(54,2).Value = 25527
$ws.Cells.Item(54,3).Value = 18
$ws.Cells.Item(54,4).Value = 23364
$ws.Cells.Item(54,5).Value = 422
$ws.Cells.Item(54,6).Value = 0
$ws.Cells.Item(54,7).Value = 0
$ws.Cells.Item(54,8).Value = 1741

$ws.Cells.Item(63,1).Value = "Argelia"
$ws.Cells.Item(63,2).Value = 15941
$ws.Cells.Item(63,3).Value = 441
$ws.Cells.Item(63,4).Value = 11492
$ws.Cells.Item(63,5).Value = 3497
$ws.Cells.Item(63,6).Value = 0
$ws.Cells.Item(63,7).Value = 6
$ws.Cells.Item(63,8).Value = 952

$ws.Cells.Item(64,1).Value = "Nepal"
$ws.Cells.Item(64,2).Value = 15784
$ws.Cells.Item(64,3).Value = 293
$ws.Cells.Item(64,4).Value = 6547
$ws.Cells.Item(64,5).Value = 9203
$ws.Cells.Item(64,6).Value = 0
$ws.Cells.Item(64,7).Value = 0
$ws.Cells.Item(64,8).Value = 34

$ws.Cells.Item(115,1).Value = "Libano"
$ws.Cells.Item(115,2).Value = 1873
$ws.Cells.Item(115,3).Value = 18
$ws.Cells.Item(115,4).Value = 1311
$ws.Cells.Item(115,5).Value = 526
$ws.Cells.Item(115,6).Value = 0
$ws.Cells.Item(115,7).Value = 1
$ws.Cells.Item(115,8).Value = 36

$ws.Cells.Item(137,1).Value = "Suazilandia"
$ws.Cells.Item(137,2).Value = 988
$ws.Cells.Item(137,3).Value = 34
$ws.Cells.Item(137,4).Value = 547
$ws.Cells.Item(137,5).Value = 428
$ws.Cells.Item(137,6).Value = 0
$ws.Cells.Item(137,7).Value = 0
$ws.Cells.Item(137,8).Value = 13

$ws.Cells.Item(138,1).Value = "Mozambique"
$ws.Cells.Item(138,2).Value = 987
$ws.Cells.Item(138,3).Value = 18
$ws.Cells.Item(138,4).Value = 256
$ws.Cells.Item(138,5).Value = 723
$ws.Cells.Item(138,6).Value = 0
$ws.Cells.Item(138,7).Value = 1
$ws.Cells.Item(138,8).Value = 8

$ws.Cells.Item(139,1).Value = "Burkina Faso"
$ws.Cells.Item(139,2).Value = 987
$ws.Cells.Item(139,3).Value = 0
$ws.Cells.Item(139,4).Value = 854
$ws.Cells.Item(139,5).Value = 80
$ws.Cells.Item(139,6).Value = 0
$ws.Cells.Item(139,7).Value = 0
$ws.Cells.Item(139,8).Value = 53

$ws.Cells.Item(140,1).Value = "Uruguay"
$ws.Cells.Item(140,2).Value = 955
$ws.Cells.Item(140,3).Value = 0
$ws.Cells.Item(140,4).Value = 840
$ws.Cells.Item(140,5).Value = 87
$ws.Cells.Item(140,6).Value = 0
$ws.Cells.Item(140,7).Value = 0
$ws.Cells.Item(140,8).Value = 28

$ws.Cells.Item(146,1).Value = "Montenegro"
$ws.Cells.Item(146,2).Value = 781
$ws.Cells.Item(146,3).Value = 61
$ws.Cells.Item(146,4).Value = 315
$ws.Cells.Item(146,5).Value = 452
$ws.Cells.Item(146,6).Value = 0
$ws.Cells.Item(146,7).Value = 1
$ws.Cells.Item(146,8).Value = 14

$ws.Cells.Item(147,1).Value = "Jamaica"
$ws.Cells.Item(147,2).Value = 728
$ws.Cells.Item(147,3).Value = 7
$ws.Cells.Item(147,4).Value = 569
$ws.Cells.Item(147,5).Value = 149
$ws.Cells.Item(147,6).Value = 0
$ws.Cells.Item(147,7).Value = 0
$ws.Cells.Item(147,8).Value = 10

$ws.Cells.Item(164,1).Value = "Birmania"
$ws.Cells.Item(164,2).Value = 313
$ws.Cells.Item(164,3).Value = 0
$ws.Cells.Item(164,4).Value = 241
$ws.Cells.Item(164,5).Value = 66
$ws.Cells.Item(164,6).Value = 0
$ws.Cells.Item(164,7).Value = 0
$ws.Cells.Item(164,8).Value = 6

$ws.Cells.Item(179,1).Value = "Trinidad yTobago"
$ws.Cells.Item(179,2).Value = 131
$ws.Cells.Item(179,3).Value = 1
$ws.Cells.Item(179,4).Value = 115
$ws.Cells.Item(179,5).Value = 8
$ws.Cells.Item(179,6).Value = 0
$ws.Cells.Item(179,7).Value = 0
$ws.Cells.Item(179,8).Value = 8

$ws.Cells.Item(186,1).Value = "Lesoto"
$ws.Cells.Item(186,2).Value = 79
$ws.Cells.Item(186,3).Value = 44
$ws.Cells.Item(186,4).Value = 11
$ws.Cells.Item(186,5).Value = 68
$ws.Cells.Item(186,6).Value = 0
$ws.Cells.Item(186,7).Value = 0
$ws.Cells.Item(186,8).Value = 0

$ws.Cells.Item(187,1).Value = "Butan"
$ws.Cells.Item(187,2).Value = 78
$ws.Cells.Item(187,3).Value = 0
$ws.Cells.Item(187,4).Value = 51
$ws.Cells.Item(187,5).Value = 27
$ws.Cells.Item(187,6).Value = 0
$ws.Cells.Item(187,7).Value = 0
$ws.Cells.Item(187,8).Value = 0

$ws.Cells.Item(188,1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(188,2).Value = 78
$ws.Cells.Item(188,3).Value = 0
$ws.Cells.Item(188,4).Value = 63
$ws.Cells.Item(188,5).Value = 0
$ws.Cells.Item(188,6).Value = 0
$ws.Cells.Item(188,7).Value = 0
$ws.Cells.Item(188,8).Value = 15

$ws.Cells.Item(205,1).Value = "Dominica"
$ws.Cells.Item(205,2).Value = 18
$ws.Cells.Item(205,3).Value = 0
$ws.Cells.Item(205,4).Value = 18
$ws.Cells.Item(205,5).Value = 0
$ws.Cells.Item(205,6).Value = 0
$ws.Cells.Item(205,7).Value = 0
$ws.Cells.Item(205,8).Value = 0

$ws.Cells.Item(206,1).Value = "Fiyi"
$ws.Cells.Item(206,2).Value = 18
$ws.Cells.Item(206,3).Value = 0
$ws.Cells.Item(206,4).Value = 18
$ws.Cells.Item(206,5).Value = 0
$ws.Cells.Item(206,6).Value = 0
$ws.Cells.Item(206,7).Value = 0
$ws.Cells.Item(206,8).Value = 0
